$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp shown in cell A1
$ws.Range("A1").Value = "Datos actualizados a 28 de Marzo de 2020 a las 10:29"

# Suiza (row 12): update totals
$ws.Range("B12").Value = 13138
$ws.Range("C12").Value = 210
$ws.Range("E12").Value = 11371
$ws.Range("G12").Value = 6
$ws.Range("H12").Value = 237

# Noruega (row 20): update totals
$ws.Range("B20").Value = 3780
$ws.Range("C20").Value = 9
$ws.Range("E20").Value = 3754

# India (row 43): update totals
$ws.Range("B43").Value = 918
$ws.Range("C43").Value = 31
$ws.Range("E43").Value = 815

# Ucrania (row 72): update totals
$ws.Range("B72").Value = 311
$ws.Range("C72").Value = 1
$ws.Range("E72").Value = 298
$ws.Range("G72").Value = 3
$ws.Range("H72").Value = 8

# Letonia (row 73): update deaths-today figure
$ws.Range("F73").Value = 3

# Taiwan now has enough cases to be placed ahead of Uruguay, Eslovaquia and
# Principado de Andorra, so rows 74-77 are rewritten: Taiwan's fresh data
# takes row 74, and the three countries that used to occupy rows 74-76 are
# pushed down to rows 75-77 respectively (their own figures are unchanged).
$ws.Range("A74").Value = "Taiwan"
$ws.Range("B74").Value = 283
$ws.Range("C74").Value = 16
$ws.Range("D74").Value = 30
$ws.Range("E74").Value = 251
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 2

$ws.Range("A75").Value = "Uruguay"
$ws.Range("B75").Value = 274
$ws.Range("C75").Value = 36
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 274
$ws.Range("F75").Value = 8
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 0

$ws.Range("A76").Value = "Eslovaquia"
$ws.Range("B76").Value = 269
$ws.Range("C76").Value = 0
$ws.Range("D76").Value = 2
$ws.Range("E76").Value = 267
$ws.Range("F76").Value = 1
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 0

$ws.Range("A77").Value = "Principado de Andorra"
$ws.Range("B77").Value = 267
$ws.Range("C77").Value = 0
$ws.Range("D77").Value = 1
$ws.Range("E77").Value = 263
$ws.Range("F77").Value = 11
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 3
